# Update the "Metadata" worksheet (sheet 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the duplicated "Contact" row (old row 11); everything below shifts up by one.
$ws.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$ws.Range("B9").Value = "Alvearie Team"

# Old "Contact" row (now row 10) becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive value (now row 14 after the row deletion above).
# Leading apostrophe forces text storage instead of Excel's auto boolean coercion.
$ws.Range("B14").Value = "'true"
